# Build site at 2021-01-29 15:15:03 UTC
#
# The course-description sheet gains a "Docentes responsáveis:" field (with
# the professor's name below it) right after "Objectives:". That pushes the
# rest of the table down two rows and, because the site generator that owns
# this sheet re-flows each field's body text into the row above it once the
# new field is spliced in, several of the long descriptive paragraphs end up
# re-attached one field earlier than before. Finally "Norma de recuperação:"
# gets a brand-new paragraph of text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert two blank rows at row 12 - everything from the old row 12
#    downward shifts to row 14 downward, carrying its row heights/formatting
#    along for the ride (matches the diff's row renumbering exactly).
$ws.Range("A12:A13").EntireRow.Insert()

# 2) Populate the two new rows with the new "responsible faculty" field.
$ws.Range("A12").Value() = "Docentes responsáveis:"
$ws.Range("B13").Value() = "5009972 - Gilberto Carvalho Coelho"
$ws.Range("C13").Value() = "5009972 - Gilberto Carvalho Coelho"

# 3) Re-flow the long body paragraphs: each one now sits one field earlier
#    than it used to (the text once attached to "Programa:" now belongs to
#    "Programa resumido:", "Método:"'s text now belongs to "Programa:", and
#    so on down the chain).
$textProgramaResumido = "A. Introdução; teoria básica de equilíbrio de fases;B. Sistemas unários;C. Sistemas binários;D. Sistemas ternários;E. Cálculo termodinâmico de diagramas de fases."
$ws.Range("B14").Value() = $textProgramaResumido
$ws.Range("C14").Value() = $textProgramaResumido

$textPrograma = "1. Introdução; revisão da termodinâmica de soluções; teoria básica de equilíbrio de fases; curvas de energia livre versus composição; regra das fases; 2. Sistemas unários, equilíbrios bi-, mono- e invariantes; 3. Sistemas binários isomorfos; a regra da alavanca; solidificação em equilíbrio e fora de equilíbrio; mínimos e máximos; 4. Sistemas eutéticos binários; solidificação e microetruturas de ligas hipoeutéticas, eutéticas e hipereutéticas; solidificação unidirecional com eutéticos; casos limites de eutéticos; 5. Sistemas eutetóides binários; solidificação e microetruturas de ligas hipoeutetóides, eutetói-des e hipereutetóides; o sistema Fe-C; 6. Sistemas monotéticos; sistemas monotetóides; sistemas metatéticos; transformações congruentes; 7. Sistemas peritéticos binários; resfriamento em equilíbrio e fora do equilíbrio de ligas peritéticas; sistemas peritetóides binários; sistemas sintéticos binários; 8. Sistemas ternários isomorfos; o triângulo de Gibbs; seções isotérmicas; projeções liquidus; seções verticais; máximos e mínimos; resfriamento em equilíbrio; 9. Equilíbrio ternário de três fases; regra da alavanca em campos trifásicos; resfriamento em equilíbrio; 10. Equilíbrio ternário de quatro fases: equilíbrio de classe I; equilíbrio de classe II e equilíbrio de classe III; 11. Transformações congruentes em sistemas ternários; sistemas ternários complexos; 12. Cálculo termodinâmico de diagramas de fases."
$ws.Range("B16").Value() = $textPrograma
$ws.Range("C16").Value() = $textPrograma

$textMetodo = "O curso será ministrado na forma de aulas expositivas e aulas práticas em laboratório envolvendo preparação de amostras e caracterização microestrutural. Os resultados das aulas práticas serão apresentados oralmente e sujeitos a avaliação (T)."
$ws.Range("B19").Value() = $textMetodo
$ws.Range("C19").Value() = $textMetodo

$textCriterio = "Serão aplicadas duas avaliações escritas (P1 e P2) que comporão a nota final (NF) juntamente com a avaliação do trabalho prático (T). O critério para a nota final é:NF=((P1*0,8)+(T*0,2)+P2*1)/2"
$ws.Range("B20").Value() = $textCriterio
$ws.Range("C20").Value() = $textCriterio

# 4) "Norma de recuperação:" gets genuinely new wording.
$textNormaRecuperacao = "Para os alunos que obtiverem 3,0≤NF<5,0, será aplicada uma avaliação de recuperação (R) que levará ao cálculo da média final (MF) com o seguinte critério:MF=(NF+R)/2"
$ws.Range("B21").Value() = $textNormaRecuperacao
$ws.Range("C21").Value() = $textNormaRecuperacao
